# Fill in the (previously empty) half-year comment / grade cells for the
# Math, English and History-of-Israel subject tables on the student
# certificate. The Torah subject table is intentionally left untouched.
#
# Each subject table has this shape (as consecutive paragraphs in document
# order):
#   [title paragraph]      e.g. "מתמטיקה"
#   [comment paragraph]    <- free text to fill in, immediately after title
#   ["ציון:" paragraph]    grade label, a couple of paragraphs later
#   [grade paragraph]      <- grade number to fill in, immediately after label
#
# Rather than relying on fixed absolute paragraph numbers, we locate each
# subject by its title text and fill the comment paragraph right after it,
# then locate the next "ציון:" label paragraph that follows the title and
# fill the grade paragraph right after that label. This keeps the script
# robust to minor structural drift while still targeting the correct cells.

$d = $word.ActiveDocument

function Find-ParagraphIndex($startIndex, $searchText) {
    $count = $d.Paragraphs.Count
    for ($p = $startIndex; $p -le $count; $p++) {
        $par = $d.Paragraphs.Item($p)
        if ($par.Range.Text.Contains($searchText)) {
            return $p
        }
    }
    return -1
}

function Set-SubjectCommentAndGrade($titleText, $commentText, $gradeText) {
    $titleIdx = Find-ParagraphIndex 1 $titleText
    if ($titleIdx -eq -1) {
        return
    }

    # Comment paragraph immediately follows the title paragraph.
    $commentIdx = $titleIdx + 1
    $d.Paragraphs.Item($commentIdx).Range.Text = $commentText

    # Grade label ("ציון:") is the next occurrence after the title; the
    # grade paragraph immediately follows that label paragraph.
    $labelIdx = Find-ParagraphIndex ($titleIdx + 1) "ציון:"
    if ($labelIdx -eq -1) {
        return
    }
    $gradeIdx = $labelIdx + 1
    $d.Paragraphs.Item($gradeIdx).Range.Text = $gradeText
}

# --- Math / מתמטיקה ---
$mathComment = "במחצית למדנו משוואות ב2 נעלמים, פיתחנו כמה שיטות לבעיה זו,בנוסף התעסקנו בבעיות תנועה וזמן ולמדנו איך להתמודד מול זאת`nהיה לנו הספקים מעולים!`nהודיה את ילדה מקסימה, שיהיה לך הרבה הצלחה בהמשך! "
Set-SubjectCommentAndGrade "מתמטיקה" $mathComment "96"

# --- English / אנגלית ---
$englishComment = "במחצית זאת התמקדנו על הבנה חזקה של הטקסטים ולמדנו את השיטות להבנת הנקרא, חזרנו על שאלות חוזרות ופיתחנו שיטות קלות לפיתרתן.`nהודיה את ילדה נהדרת, הרבה הצלחה!"
Set-SubjectCommentAndGrade "אנגלית" $englishComment "86"

# --- History of Israel / תולדות ישראל ---
$historyComment = "במחצית זאת למדנו על גדולי ישראל בכל מיני יבשות, על המצב של היהודים בתקופות שלטון שונות,`nהודיה הרבה הצלחה!"
Set-SubjectCommentAndGrade "תולדות ישראל" $historyComment "98"
